$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new product row before row 32 (shifts the existing
# "سرنجات 3 سم" / "شفرات فينوس حريمي " / "مناديل سولو سحب صغيره" rows,
# the totals row and the footer row down by one).
$ws.Range("A32:Q32").Insert(-4121)

# Bring over the same cell formatting used by the other product rows
# (copy from what is now row 33, formerly row 32) so the new row matches
# the existing style ids instead of picking up brand-new ones.
$ws.Range("A33:Q33").Copy()
$ws.Range("A32:Q32").PasteSpecial(-4122)

# Recreate the merged cells for the new row, matching the layout used by
# every other product row.
$ws.Range("A32:B32").Merge()
$ws.Range("C32:G32").Merge()
$ws.Range("H32:K32").Merge()
$ws.Range("L32:M32").Merge()
$ws.Range("N32:O32").Merge()

# Match the row height used by the other product rows.
$ws.Rows.Item(32).RowHeight = 25.5

# Populate the new "VOLTAREN 75MG/3ML 3 AMP." row.
$ws.Cells.Item(32, 1).Value = 26
$ws.Cells.Item(32, 3).Value = "VOLTAREN 75MG/3ML 3 AMP."
$ws.Cells.Item(32, 8).Value = "3:3"
$ws.Cells.Item(32, 12).Value = "1"
$ws.Cells.Item(32, 14).Value = "51.00"
$ws.Cells.Item(32, 16).Value = "16.8300"
$ws.Cells.Item(32, 17).Value = "0:1"

# Update the "سرنجات 3 سم" row (now shifted down to row 33) with the new
# sell price / transaction-count figures.
$ws.Cells.Item(33, 16).Value = "16.0000"
$ws.Cells.Item(33, 17).Value = "8:0"

# Update the totals row (now shifted down to row 36) with the new sum.
$ws.Cells.Item(36, 16).Value = 2360.78

# Update the generated timestamp string in the footer (now shifted down
# to row 37 along with the rest of the footer row).
$ws.Cells.Item(37, 1).Value = "Sunday, 28 September, 2025 3:29 PM"
